# Refresh the event-tracker snapshot: bump "want to go" counts (column F)
# across the four sheets, and drop the no-longer-relevant LoveLive live
# stream row (it was "不可售" / not-for-sale) from 本地生活 and from the
# combined 全部类型 view, shifting everything below it up by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions) - F-column ("想去人数") value refresh
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 326
$ws1.Cells.Item(3, 6).Value = 277
$ws1.Cells.Item(5, 6).Value = 174
$ws1.Cells.Item(6, 6).Value = 664
$ws1.Cells.Item(8, 6).Value = 476
$ws1.Cells.Item(9, 6).Value = 81
$ws1.Cells.Item(10, 6).Value = 518
$ws1.Cells.Item(11, 6).Value = 392
$ws1.Cells.Item(12, 6).Value = 67
$ws1.Cells.Item(14, 6).Value = 112
$ws1.Cells.Item(15, 6).Value = 198

# ---------------------------------------------------------------
# Sheet 2: 演出 (Performances) - F-column value refresh
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(8, 6).Value = 41
$ws2.Cells.Item(9, 6).Value = 182
$ws2.Cells.Item(14, 6).Value = 28
$ws2.Cells.Item(15, 6).Value = 40

# ---------------------------------------------------------------
# Sheet 3: 本地生活 (Local life) - F-column value refresh, then drop
# row 6 (2024-10-06 LoveLive live stream, marked 不可售)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 6222
$ws3.Cells.Item(4, 6).Value = 755
$ws3.Cells.Item(5, 6).Value = 1819
$ws3.Rows.Item(6).Delete()

# ---------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - same LoveLive row removal, then the
# same F-column refresh applied at the post-deletion row positions
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 6222
$ws4.Cells.Item(4, 6).Value = 755
$ws4.Cells.Item(5, 6).Value = 1819
$ws4.Rows.Item(6).Delete()

# The "#" index column (A) in this sheet is a fixed row-position counter
# (row number - 1), not data tied to a particular event; it is left as-is
# by the source refresh, so put it back after the native row-shift moved
# it along with everything else.
for ($r = 6; $r -le 35; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

$ws4.Cells.Item(6, 6).Value = 326
$ws4.Cells.Item(7, 6).Value = 277
$ws4.Cells.Item(12, 6).Value = 174
$ws4.Cells.Item(15, 6).Value = 664
$ws4.Cells.Item(18, 6).Value = 41
$ws4.Cells.Item(19, 6).Value = 476
$ws4.Cells.Item(20, 6).Value = 182
$ws4.Cells.Item(21, 6).Value = 81
$ws4.Cells.Item(22, 6).Value = 518
$ws4.Cells.Item(24, 6).Value = 392
$ws4.Cells.Item(25, 6).Value = 67
$ws4.Cells.Item(29, 6).Value = 112
$ws4.Cells.Item(31, 6).Value = 28
$ws4.Cells.Item(32, 6).Value = 40
$ws4.Cells.Item(35, 6).Value = 198
